$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1673469387755102
$ws.Range("C2").Value = 0.6489795918367347
$ws.Range("P2").Value = 0.1061224489795918
$ws.Range("S2").Value = 0.07755102040816327

# Row 3
$ws.Range("B3").Value = 0.006172839506172839
$ws.Range("C3").Value = 0.04320987654320987
$ws.Range("J3").Value = 0.006172839506172839
$ws.Range("P3").Value = 0.7654320987654321
$ws.Range("S3").Value = 0.1790123456790123

# Row 4
$ws.Range("P4").Value = 0.6875
$ws.Range("S4").Value = 0.3125

# Row 6
$ws.Range("B6").Value = 0.06637168141592921
$ws.Range("D6").Value = 0.008849557522123894
$ws.Range("F6").Value = 0.03539823008849557
$ws.Range("J6").Value = 0.2256637168141593
$ws.Range("O6").Value = 0.02212389380530973
$ws.Range("Q6").Value = 0.1814159292035398
$ws.Range("R6").Value = 0.05309734513274336
$ws.Range("S6").Value = 0.4070796460176991

# Row 7
$ws.Range("B7").Value = 0.09782608695652174
$ws.Range("D7").Value = 0.01630434782608696
$ws.Range("F7").Value = 0.04891304347826087
$ws.Range("J7").Value = 0.09239130434782608
$ws.Range("O7").Value = 0.01630434782608696
$ws.Range("Q7").Value = 0.2608695652173913
$ws.Range("R7").Value = 0.1141304347826087
$ws.Range("S7").Value = 0.3532608695652174

# Row 8
$ws.Range("B8").Value = 0.07434052757793765
$ws.Range("D8").Value = 0.01199040767386091
$ws.Range("F8").Value = 0.07913669064748201
$ws.Range("J8").Value = 0.105515587529976
$ws.Range("O8").Value = 0.007194244604316547
$ws.Range("Q8").Value = 0.1558752997601918
$ws.Range("R8").Value = 0.1007194244604317
$ws.Range("S8").Value = 0.4652278177458034

# Row 9
$ws.Range("B9").Value = 0.09363295880149813
$ws.Range("D9").Value = 0.02247191011235955
$ws.Range("E9").Value = 0.003745318352059925
$ws.Range("F9").Value = 0.0749063670411985
$ws.Range("J9").Value = 0.101123595505618
$ws.Range("O9").Value = 0.0299625468164794
$ws.Range("Q9").Value = 0.1535580524344569
$ws.Range("R9").Value = 0.09363295880149813
$ws.Range("S9").Value = 0.4269662921348314

# Row 10
$ws.Range("B10").Value = 0.09722222222222222
$ws.Range("D10").Value = 0.01909722222222222
$ws.Range("F10").Value = 0.07725694444444445
$ws.Range("J10").Value = 0.09722222222222222
$ws.Range("O10").Value = 0.01302083333333333
$ws.Range("Q10").Value = 0.2161458333333333
$ws.Range("R10").Value = 0.07378472222222222
$ws.Range("S10").Value = 0.40625

# Row 11
$ws.Range("G11").Value = 0.13
$ws.Range("J11").Value = 0.1266666666666667
$ws.Range("K11").Value = 0.1933333333333333
$ws.Range("L11").Value = 0.5366666666666666
$ws.Range("S11").Value = 0.01333333333333333

# Row 12
$ws.Range("G12").Value = 0.7011494252873564
$ws.Range("J12").Value = 0.2183908045977012
$ws.Range("K12").Value = 0.01724137931034483
$ws.Range("L12").Value = 0.03448275862068965
$ws.Range("S12").Value = 0.02873563218390805

# Row 15
$ws.Range("F15").Value = 0.03043478260869565
$ws.Range("H15").Value = 0.1695652173913043
$ws.Range("I15").Value = 0.1347826086956522
$ws.Range("J15").Value = 0.3130434782608696
$ws.Range("K15").Value = 0.07391304347826087
$ws.Range("M15").Value = 0.008695652173913044
$ws.Range("O15").Value = 0.04347826086956522
$ws.Range("S15").Value = 0.2260869565217391

# Row 16
$ws.Range("F16").Value = 0.01197604790419162
$ws.Range("H16").Value = 0.125748502994012
$ws.Range("I16").Value = 0.1077844311377246
$ws.Range("J16").Value = 0.407185628742515
$ws.Range("K16").Value = 0.08982035928143713
$ws.Range("M16").Value = 0.01197604790419162
$ws.Range("O16").Value = 0.08383233532934131
$ws.Range("S16").Value = 0.1616766467065868

# Row 17
$ws.Range("F17").Value = 0.02727272727272727
$ws.Range("H17").Value = 0.1636363636363636
$ws.Range("I17").Value = 0.1318181818181818
$ws.Range("J17").Value = 0.3772727272727273
$ws.Range("K17").Value = 0.09545454545454546
$ws.Range("M17").Value = 0.02045454545454545
$ws.Range("O17").Value = 0.07272727272727272
$ws.Range("S17").Value = 0.1113636363636364

# Row 18
$ws.Range("F18").Value = 0.01648351648351648
$ws.Range("H18").Value = 0.1758241758241758
$ws.Range("I18").Value = 0.09340659340659341
$ws.Range("J18").Value = 0.4010989010989011
$ws.Range("K18").Value = 0.08791208791208792
$ws.Range("M18").Value = 0.01098901098901099
$ws.Range("O18").Value = 0.05494505494505494
$ws.Range("S18").Value = 0.1593406593406593

# Row 19
$ws.Range("F19").Value = 0.01660377358490566
$ws.Range("H19").Value = 0.190188679245283
$ws.Range("I19").Value = 0.1071698113207547
$ws.Range("J19").Value = 0.3403773584905661
$ws.Range("K19").Value = 0.1109433962264151
$ws.Range("M19").Value = 0.01660377358490566
$ws.Range("N19").Value = 0.002264150943396227
$ws.Range("O19").Value = 0.07698113207547169
$ws.Range("S19").Value = 0.1388679245283019

Write-Output "Applied team specific time commit updates"